$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Moyenne de l'etudiant" (column E) values for the student rows.
# Mapping: row number -> new average value, as per the commit's grade corrections.
$ws.Cells.Item(3, 5).Value = 20
$ws.Cells.Item(4, 5).Value = 14
$ws.Cells.Item(5, 5).Value = 19
$ws.Cells.Item(6, 5).Value = 6
$ws.Cells.Item(7, 5).Value = 16
$ws.Cells.Item(8, 5).Value = 11
$ws.Cells.Item(9, 5).Value = 14
$ws.Cells.Item(10, 5).Value = 12
$ws.Cells.Item(11, 5).Value = 6
$ws.Cells.Item(12, 5).Value = 10
$ws.Cells.Item(13, 5).Value = 5
$ws.Cells.Item(14, 5).Value = 7
$ws.Cells.Item(15, 5).Value = 9
$ws.Cells.Item(16, 5).Value = 15
$ws.Cells.Item(17, 5).Value = 10
$ws.Cells.Item(18, 5).Value = 7
$ws.Cells.Item(19, 5).Value = 8
$ws.Cells.Item(20, 5).Value = 9
$ws.Cells.Item(21, 5).Value = 10
$ws.Cells.Item(22, 5).Value = 6
$ws.Cells.Item(23, 5).Value = 7
$ws.Cells.Item(24, 5).Value = 20
$ws.Cells.Item(25, 5).Value = 18
$ws.Cells.Item(27, 5).Value = 20
$ws.Cells.Item(28, 5).Value = 8
$ws.Cells.Item(30, 5).Value = 17
$ws.Cells.Item(31, 5).Value = 8
$ws.Cells.Item(32, 5).Value = 8
$ws.Cells.Item(33, 5).Value = 8
$ws.Cells.Item(34, 5).Value = 14
$ws.Cells.Item(35, 5).Value = 19
$ws.Cells.Item(36, 5).Value = 13
$ws.Cells.Item(37, 5).Value = 5
$ws.Cells.Item(38, 5).Value = 15
$ws.Cells.Item(39, 5).Value = 19
$ws.Cells.Item(41, 5).Value = 13
$ws.Cells.Item(42, 5).Value = 16
$ws.Cells.Item(43, 5).Value = 5
$ws.Cells.Item(44, 5).Value = 14
$ws.Cells.Item(45, 5).Value = 16
$ws.Cells.Item(46, 5).Value = 16
$ws.Cells.Item(47, 5).Value = 15
$ws.Cells.Item(48, 5).Value = 11
$ws.Cells.Item(49, 5).Value = 19
$ws.Cells.Item(50, 5).Value = 11
$ws.Cells.Item(52, 5).Value = 5
$ws.Cells.Item(53, 5).Value = 10
$ws.Cells.Item(54, 5).Value = 20
$ws.Cells.Item(55, 5).Value = 5
$ws.Cells.Item(56, 5).Value = 20
$ws.Cells.Item(57, 5).Value = 15
$ws.Cells.Item(58, 5).Value = 17
$ws.Cells.Item(59, 5).Value = 16
$ws.Cells.Item(60, 5).Value = 18
$ws.Cells.Item(61, 5).Value = 5
$ws.Cells.Item(62, 5).Value = 5
$ws.Cells.Item(63, 5).Value = 8
